$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: B2 and F2 were auto-generated pandas placeholders
# ("unnamed: 1_level_1" / "unnamed: 5_level_1"); they should read "total",
# same label as C2.
$totalLabel = $ws.Range("C2").Value2
$ws.Range("B2").Value = $totalLabel
$ws.Range("F2").Value = $totalLabel

# Remove the two blank-data rows (originally rows 5 and 8) so the real
# data rows that followed them shift up and the table becomes contiguous.
# Row 8 is deleted first so row 5's index is still valid afterwards.
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(5).Delete()
